$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Do you feel that this access control request is actually a vulnerability?"
#    The two grammar-flagged runs ("...actually " / "a vulnerability") get
#    merged into a single run and the wrapping grammar-check markers
#    (<w:proofErr w:type="gramStart"/> / <w:proofErr w:type="gramEnd"/>) are
#    cleared, as if the user accepted/ignored the grammar squiggle while
#    touching up the sentence. The trailing "?" run is left alone.
# ---------------------------------------------------------------------------
$q = "Do you feel that this access control request is actually a vulnerability?"
$d.Content.Find.Execute($q, $false, $false, $false, $false, $false, $true, 1, $false, $q, 2)

# Re-isolate the trailing "?" back into its own run (it was merged above along
# with the rest of the sentence while the grammar markers were being swept
# away) by forcing a run split right before it via a harmless formatting
# round-trip.
$rng = $d.Content
$rng.Find.Execute($q, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$qMark = $d.Range($rng.End - 1, $rng.End)
$qMark.Bold = 1
$qMark.Bold = 0

# ---------------------------------------------------------------------------
# 2) Drop the stray trailing-space run after "...attached to this line of
#    code?" in the warning section (leaves the sentence's own run intact).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Why do you think this warning is attached to this line of code? ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$trailingSpace = $d.Range($rng2.End - 1, $rng2.End)
$trailingSpace.Delete()

# ---------------------------------------------------------------------------
# 3) Add a new closing interview question after "Do you think you would use
#    ASIDE in the real world?".
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Do you think you would use ASIDE in the real world?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3.InsertParagraphAfter()
$newPara = $rng3.Paragraphs(1).Next()
$newPara.Range.Text = "Is there anything else you want to say about your experience with ASIDE today?"
